# Fill in the "04dec2025" column (E) counts for several empadronadores that
# had not been entered yet (they were all placeholder 0's).
#
# The workbook carries the same numbers twice:
#   - "crosstab" stores them as real numbers.
#   - "annot"    stores them as text (it renders 0 as a blank cell), so the
#                same figures need to be written there too, but as strings.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 1
    6  = 14
    7  = 13
    8  = 10
    9  = 10
    10 = 10
    12 = 14
}

$crosstab = $wb.Worksheets.Item("crosstab")
$annot    = $wb.Worksheets.Item("annot")

foreach ($row in $updates.Keys) {
    $value = $updates[$row]

    # Numeric sheet: plain value assignment.
    $crosstab.Range("E$row").Value = $value

    # Text sheet: force text formatting first so the numeric-looking string
    # isn't re-interpreted as a number, then restore the default "Normal"
    # style so no stray formatting is left behind on the cell.
    $cell = $annot.Range("E$row")
    $cell.NumberFormat = "@"
    $cell.Value = "$value"
    $cell.Style = "Normal"
}
